$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mechanism")

# --- Panel row-label text edits -----------------------------------------
# "Forced commitment" -> "Mandatory structured" (rows 5 and 17, column A)
$ws.Range("A5").Value = "Mandatory structured"
$ws.Range("A17").Value = "Mandatory structured"

# "Choice commitment" -> "Choice " (rows 7 and 19, column A)
$ws.Range("A7").Value = "Choice "
$ws.Range("A19").Value = "Choice "

# Row 18 column A used to hold the (now unused) empty-string shared entry;
# the cell is removed entirely in the target sheet.
$ws.Range("A18").Value = ""

# --- Selection / view state ----------------------------------------------
$ws.Activate()
$ws.Range("A2:H24").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
